$d = $word.ActiveDocument

$r = $d.Content
$ok = $r.Find.Execute("remover as tarefas.", $false, $false, $false, $false, $false, $true, 1, $false, "remover e alterar as tarefas. Além disso o sistema conta com uma tela de ajuda.", 2)
Write-Output $ok
